# Fruta / hortaliza, semanal
# Insert 3 new weekly rows (370-372) into the daily logic subconjunto sheet
# for "Terminal La Palmera de La Serena - Pera", pushing the existing rows
# 370-383 down to 373-386.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows above the current row 370 (shifts 370:383 -> 373:386)
$ws.Rows.Item(370).Insert()
$ws.Rows.Item(371).Insert()
$ws.Rows.Item(372).Insert()

# --- New row 370 ---
$ws.Cells.Item(370, 1).Value = 8
$ws.Cells.Item(370, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(370, 3).Value = "Coquimbo"
$ws.Cells.Item(370, 4).Value = 44516
$ws.Cells.Item(370, 5).Value = 4
$ws.Cells.Item(370, 6).Value = "Fruta"
$ws.Cells.Item(370, 7).Value = 100104
$ws.Cells.Item(370, 8).Value = "Frutos de pepita"
$ws.Cells.Item(370, 9).Value = 100104005
$ws.Cells.Item(370, 10).Value = "Pera"
$ws.Cells.Item(370, 11).Value = "Packham's Triumph"
$ws.Cells.Item(370, 12).Value = "Especial"
$ws.Cells.Item(370, 13).Value = 20
$ws.Cells.Item(370, 14).Value = 290000
$ws.Cells.Item(370, 15).Value = 300000
$ws.Cells.Item(370, 16).Value = 295000
$ws.Cells.Item(370, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(370, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(370, 19).Value = 656
$ws.Cells.Item(370, 20).Value = 450

# --- New row 371 ---
$ws.Cells.Item(371, 1).Value = 8
$ws.Cells.Item(371, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(371, 3).Value = "Coquimbo"
$ws.Cells.Item(371, 4).Value = 44516
$ws.Cells.Item(371, 5).Value = 4
$ws.Cells.Item(371, 6).Value = "Fruta"
$ws.Cells.Item(371, 7).Value = 100104
$ws.Cells.Item(371, 8).Value = "Frutos de pepita"
$ws.Cells.Item(371, 9).Value = 100104005
$ws.Cells.Item(371, 10).Value = "Pera"
$ws.Cells.Item(371, 11).Value = "Packham's Triumph"
$ws.Cells.Item(371, 12).Value = "Primera"
$ws.Cells.Item(371, 13).Value = 16
$ws.Cells.Item(371, 14).Value = 270000
$ws.Cells.Item(371, 15).Value = 280000
$ws.Cells.Item(371, 16).Value = 275000
$ws.Cells.Item(371, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(371, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(371, 19).Value = 611
$ws.Cells.Item(371, 20).Value = 450

# --- New row 372 ---
$ws.Cells.Item(372, 1).Value = 8
$ws.Cells.Item(372, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(372, 3).Value = "Coquimbo"
$ws.Cells.Item(372, 4).Value = 44516
$ws.Cells.Item(372, 5).Value = 4
$ws.Cells.Item(372, 6).Value = "Fruta"
$ws.Cells.Item(372, 7).Value = 100104
$ws.Cells.Item(372, 8).Value = "Frutos de pepita"
$ws.Cells.Item(372, 9).Value = 100104005
$ws.Cells.Item(372, 10).Value = "Pera"
$ws.Cells.Item(372, 11).Value = "Packham's Triumph"
$ws.Cells.Item(372, 12).Value = "Segunda"
$ws.Cells.Item(372, 13).Value = 14
$ws.Cells.Item(372, 14).Value = 250000
$ws.Cells.Item(372, 15).Value = 260000
$ws.Cells.Item(372, 16).Value = 255000
$ws.Cells.Item(372, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(372, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(372, 19).Value = 567
$ws.Cells.Item(372, 20).Value = 450
